# Auto-generated edit script applying the Marilith_Profits.xlsx diff
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H-N) across all 8 sheets
$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 629.1667
$ws.Range("J29").Value = 3500
$ws.Range("L29").Value = 10500
$ws.Range("N29").Value = -11062
$ws.Range("H32").Value = 946.3158
$ws.Range("I32").Value = 798.93335
$ws.Range("J32").Value = 1499
$ws.Range("K32").Value = 798.93335
$ws.Range("L32").Value = 1499
$ws.Range("M32").Value = -472.93335
$ws.Range("N32").Value = -2151
$ws.Range("H62").Value = 6183.75
$ws.Range("J62").Value = 6557
$ws.Range("L62").Value = 6557
$ws.Range("N62").Value = -7805
$ws.Range("H65").Value = 6183.75
$ws.Range("J65").Value = 6557
$ws.Range("L65").Value = 32785
$ws.Range("N65").Value = -39025
$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = ""
$ws.Range("N69").Value = 0
$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = ""
$ws.Range("N72").Value = 0
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = ""
$ws.Range("M74").Value = ""
$ws.Range("N74").Value = 0
$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = ""
$ws.Range("M77").Value = ""
$ws.Range("N77").Value = 0
$ws.Range("H111").Value = 530.8333
$ws.Range("I111").Value = 451.66666
$ws.Range("J111").Value = 768.3333
$ws.Range("K111").Value = 1354.99998
$ws.Range("L111").Value = 2304.9999
$ws.Range("M111").Value = 1712.00002
$ws.Range("N111").Value = -8438.999899999999
$ws.Range("H113").Value = 4089.9546
$ws.Range("I113").Value = 2698.6
$ws.Range("J113").Value = 7071.4287
$ws.Range("K113").Value = 2698.6
$ws.Range("L113").Value = 7071.4287
$ws.Range("M113").Value = 555.4000000000001
$ws.Range("N113").Value = -13579.4287
$ws.Range("H137").Value = 4517.6665
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 4517.6665
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = ""
$ws.Range("M137").Value = 13552.9995
$ws.Range("N137").Value = -18652.9995

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H33").Value = 7513.75
$ws.Range("J33").Value = 25029
$ws.Range("L33").Value = 25029
$ws.Range("N33").Value = -25687
$ws.Range("H88").Value = 1982.1578
$ws.Range("I88").Value = 380.125
$ws.Range("J88").Value = 3147.2727
$ws.Range("K88").Value = 380.125
$ws.Range("L88").Value = 3147.2727
$ws.Range("M88").Value = 25.875
$ws.Range("N88").Value = -3959.2727
$ws.Range("H91").Value = 1982.1578
$ws.Range("I91").Value = 380.125
$ws.Range("J91").Value = 3147.2727
$ws.Range("K91").Value = 380.125
$ws.Range("L91").Value = 3147.2727
$ws.Range("M91").Value = 1023.875
$ws.Range("N91").Value = -5955.2727
$ws.Range("H110").Value = 3089.75
$ws.Range("I110").Value = 2497.2
$ws.Range("K110").Value = 2497.2
$ws.Range("M110").Value = -452.1999999999998
$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = ""
$ws.Range("N118").Value = 0
$ws.Range("H132").Value = 4240.3335
$ws.Range("I132").Value = 3716.7273
$ws.Range("J132").Value = 10000
$ws.Range("K132").Value = 11150.1819
$ws.Range("L132").Value = 30000
$ws.Range("M132").Value = -8620.1819
$ws.Range("N132").Value = -35060

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H15").Value = 33736.25
$ws.Range("J15").Value = 35984.285
$ws.Range("L15").Value = 35984.285
$ws.Range("N15").Value = -36438.285
$ws.Range("H19").Value = 25979.166
$ws.Range("J19").Value = 25979.166
$ws.Range("L19").Value = 25979.166
$ws.Range("N19").Value = -26325.166
$ws.Range("H105").Value = 1943.7368
$ws.Range("I105").Value = 1923.9445
$ws.Range("K105").Value = 1923.9445
$ws.Range("M105").Value = -176.9445000000001
$ws.Range("H107").Value = 1288.8334
$ws.Range("I107").Value = 1288.8334
$ws.Range("K107").Value = 1288.8334
$ws.Range("M107").Value = 631.1666
$ws.Range("H134").Value = 6593.2
$ws.Range("I134").Value = 6349.857
$ws.Range("J134").Value = 10000
$ws.Range("K134").Value = 19049.571
$ws.Range("L134").Value = 30000
$ws.Range("M134").Value = -16514.571
$ws.Range("N134").Value = -35070

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = ""
$ws.Range("N18").Value = 0
$ws.Range("H56").Value = 13333
$ws.Range("I56").Value = 7499.5
$ws.Range("K56").Value = 7499.5
$ws.Range("M56").Value = -6654.5
$ws.Range("H62").Value = 1799.25
$ws.Range("I62").Value = 1799.25
$ws.Range("K62").Value = 1799.25
$ws.Range("M62").Value = -1175.25
$ws.Range("H65").Value = 1799.25
$ws.Range("I65").Value = 1799.25
$ws.Range("K65").Value = 8996.25
$ws.Range("M65").Value = -5876.25
$ws.Range("H99").Value = 3799.625
$ws.Range("I99").Value = 3485.2856
$ws.Range("J99").Value = 6000
$ws.Range("K99").Value = 3485.2856
$ws.Range("L99").Value = 6000
$ws.Range("M99").Value = -1987.2856
$ws.Range("N99").Value = -8996
$ws.Range("H107").Value = 779
$ws.Range("I107").Value = 753.5333000000001
$ws.Range("J107").Value = 874.5
$ws.Range("K107").Value = 753.5333000000001
$ws.Range("L107").Value = 874.5
$ws.Range("M107").Value = 1166.4667
$ws.Range("N107").Value = -4714.5
$ws.Range("H126").Value = 3799.625
$ws.Range("I126").Value = 3485.2856
$ws.Range("J126").Value = 6000
$ws.Range("K126").Value = 10455.8568
$ws.Range("L126").Value = 18000
$ws.Range("M126").Value = -7985.856800000001
$ws.Range("N126").Value = -22940
$ws.Range("H132").Value = 3629.1052
$ws.Range("I132").Value = 3463.6
$ws.Range("K132").Value = 10390.8
$ws.Range("M132").Value = -7860.799999999999

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 44.153847
$ws.Range("J12").Value = 64.57143000000001
$ws.Range("L12").Value = 193.71429
$ws.Range("N12").Value = -539.71429
$ws.Range("H120").Value = 949
$ws.Range("I120").Value = 949
$ws.Range("K120").Value = 2847
$ws.Range("M120").Value = 1991
$ws.Range("H132").Value = 3282.1428
$ws.Range("I132").Value = 3900
$ws.Range("J132").Value = 3179.1667
$ws.Range("K132").Value = 35100
$ws.Range("L132").Value = 28612.5003
$ws.Range("M132").Value = -32570
$ws.Range("N132").Value = -33672.5003

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H35").Value = 3025
$ws.Range("I35").Value = 3025
$ws.Range("K35").Value = 3025
$ws.Range("M35").Value = -2727
$ws.Range("H107").Value = 166.8
$ws.Range("I107").Value = 166.8
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 166.8
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = ""
$ws.Range("N107").Value = 1753.2

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 17655.625
$ws.Range("I7").Value = 18857.5
$ws.Range("K7").Value = 18857.5
$ws.Range("M7").Value = -18745.5
$ws.Range("H22").Value = 753.3103599999999
$ws.Range("I22").Value = 666.6667
$ws.Range("J22").Value = 814.4706
$ws.Range("K22").Value = 666.6667
$ws.Range("L22").Value = 814.4706
$ws.Range("M22").Value = -371.6667
$ws.Range("N22").Value = -1404.4706
$ws.Range("H27").Value = 753.3103599999999
$ws.Range("I27").Value = 666.6667
$ws.Range("J27").Value = 814.4706
$ws.Range("K27").Value = 666.6667
$ws.Range("L27").Value = 814.4706
$ws.Range("M27").Value = -559.6667
$ws.Range("N27").Value = -1028.4706
$ws.Range("H55").Value = 278
$ws.Range("J55").Value = 278
$ws.Range("L55").Value = 278
$ws.Range("N55").Value = -624
$ws.Range("H61").Value = 2117.6667
$ws.Range("I61").Value = 2392.5715
$ws.Range("J61").Value = 1155.5
$ws.Range("K61").Value = 2392.5715
$ws.Range("L61").Value = 1155.5
$ws.Range("M61").Value = -2190.5715
$ws.Range("N61").Value = -1559.5
$ws.Range("H82").Value = 744.8333
$ws.Range("I82").Value = 594
$ws.Range("K82").Value = 594
$ws.Range("M82").Value = -233
$ws.Range("H85").Value = 744.8333
$ws.Range("I85").Value = 594
$ws.Range("K85").Value = 594
$ws.Range("M85").Value = 654
$ws.Range("H93").Value = 975
$ws.Range("I93").Value = 975
$ws.Range("K93").Value = 975
$ws.Range("M93").Value = 273
$ws.Range("H113").Value = 2117.6667
$ws.Range("I113").Value = 2392.5715
$ws.Range("J113").Value = 1155.5
$ws.Range("K113").Value = 2392.5715
$ws.Range("L113").Value = 1155.5
$ws.Range("M113").Value = -222.5715
$ws.Range("N113").Value = -5495.5
$ws.Range("H126").Value = 17655.625
$ws.Range("I126").Value = 18857.5
$ws.Range("K126").Value = 56572.5
$ws.Range("M126").Value = -54102.5

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 1020540.8
$ws.Range("I4").Value = 1250751
$ws.Range("J4").Value = 99700
$ws.Range("K4").Value = 1250751
$ws.Range("L4").Value = 99700
$ws.Range("M4").Value = -1250638
$ws.Range("N4").Value = -99926
